$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in F1, copying the same style as the other
# header cells (bold font, border, centered alignment) from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F23).
$timestamps = @(
    "2021-10-05 13:39:09.965085",
    "2021-10-05 13:39:09.965096",
    "2021-10-05 13:39:09.965100",
    "2021-10-05 13:39:09.965102",
    "2021-10-05 13:39:09.965105",
    "2021-10-05 13:39:09.965108",
    "2021-10-05 13:39:09.965111",
    "2021-10-05 13:39:09.965114",
    "2021-10-05 13:39:09.965117",
    "2021-10-05 13:39:09.965119",
    "2021-10-05 13:39:09.965122",
    "2021-10-05 13:39:09.965124",
    "2021-10-05 13:39:09.965127",
    "2021-10-05 13:39:09.965129",
    "2021-10-05 13:39:09.965132",
    "2021-10-05 13:39:09.965134",
    "2021-10-05 13:39:09.965137",
    "2021-10-05 13:39:09.965140",
    "2021-10-05 13:39:09.965143",
    "2021-10-05 13:39:09.965146",
    "2021-10-05 13:39:09.965148",
    "2021-10-05 13:39:09.965151"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
